$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

# Citation count update
Replace-Text "cited 1203 times" "cited 1236 times"

# h-index update
Replace-Text "-index is 16 and my" "-index is 17 and my"

# i10-index update
Replace-Text "-index is 25" "-index is 26"

# Updated date
Replace-Text "(Updated 01.06.2021)" "(Updated 19.06.2021)"

# Bibliography entry title-case -> sentence-case fixes
Replace-Text "Hva Betyr Det for Barn å Vokse Opp i" "Hva betyr det for barn å vokse opp i"
Replace-Text "En Fattig Familie? [What Does It Mean for Children to Grow up in a Poor" "en fattig familie? [What does it mean for children to grow up in a poor"
Replace-Text "Family?]" "family?]"
Replace-Text "Helsetasjonstjenesten: Barns Psykiske Helse Og" "Helsetasjonstjenesten: Barns psykiske helse og"
Replace-Text "Utvikling" "utvikling"
